$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 750; this shifts the existing rows 750-834
# down to 751-835 (and the used dimension grows from R834 to R835).
$ws.Rows.Item(750).Insert()

# Populate the newly inserted row 750 with the new weekly data point.
$ws.Cells.Item(750, 1).Value = 8
$ws.Cells.Item(750, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(750, 3).Value = "Coquimbo"
$ws.Cells.Item(750, 4).Value = 44918
$ws.Cells.Item(750, 5).Value = 4
$ws.Cells.Item(750, 6).Value = 100112045
$ws.Cells.Item(750, 7).Value = "Zapallo"
$ws.Cells.Item(750, 8).Value = "Camote"
$ws.Cells.Item(750, 9).Value = "1a nueva(o)"
$ws.Cells.Item(750, 10).Value = 2000
$ws.Cells.Item(750, 11).Value = 900
$ws.Cells.Item(750, 12).Value = 1000
$ws.Cells.Item(750, 13).Value = 950
$ws.Cells.Item(750, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(750, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(750, 16).Value = 950
$ws.Cells.Item(750, 17).Value = 1
$ws.Cells.Item(750, 18).Value = "Hortaliza"
